# "merge sa trimiti email" - update the Email column for the two users
# so the merge/mail step sends to the correct addresses, and leave the
# selection on the cell the author ended up clicking (F4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "andrei@mail.com"
$ws.Range("D3").Value = "olivia@mail.com"

$ws.Range("F4").Select()
